# daily auto push: 2026-02-14 07:01 UTC
#
# A new hourly observation for 2026/02/14 (Saturday, hour bucket "14")
# was appended to the data feed. Because the sheet is already sorted by
# date/time, the new record lands in the middle of the table (row 818),
# pushing every row below it down by one and extending the sheet's
# dimension from D859 to D860.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 818; rows 818..859 shift down to 819..860.
$ws.Rows(818).Insert()

# Fill in the new row. The date column holds plain text like "2026/02/14"
# (not a real Excel date), so a leading apostrophe keeps Excel from
# reinterpreting it as a date serial; resetting the style back to Normal
# afterwards clears the "quote prefix" flag that the apostrophe leaves
# behind, matching the plain, unstyled cells used elsewhere in column A.
$ws.Cells.Item(818, 1).Value = "'2026/02/14"
$ws.Cells.Item(818, 1).Style = "Normal"
$ws.Cells.Item(818, 2).Value = "土"
$ws.Cells.Item(818, 3).Value = 14
$ws.Cells.Item(818, 4).Value = 201
